# Log.xlsx maintenance update:
#  - record the 2/2/2015 date against the two "Receive assignments" rows
#  - widen the Date column slightly to fit
#  - log a new entry: Shaurya / Send Project Ara info (2/26/2015)
#  - leave the selection on the next empty row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: "Receive assignments" (Andrew) now has its date filled in ---
$ws.Range("B3").Value = 42037

# --- Row 4: "Receive assignments" (BOTH) gets the same date. This cell is
#     brand new, so first clone B3's formatting (date number format / style)
#     via a formats-only paste, then write the value. ---
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = 42037

# --- New row 10: Shaurya / Send Project Ara info, dated 2/26/2015 ---
$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = 42061
$ws.Range("C10").Value = "Shaurya"
$ws.Range("D10").Value = "Send Project Ara info"

# --- Column B (the Date column) is slightly wider now ---
$ws.Columns.Item(2).ColumnWidth = 14

# --- Move the selection down to the next blank row, ready for new input ---
$ws.Range("B11").Select() | Out-Null
